$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.068555616624326
$ws.Range("D2").Value = 1.071417616150721
$ws.Range("E2").Value = 1.081300175554434
$ws.Range("F2").Value = 1.086502401336667
$ws.Range("I2").Value = 1.058329338103405
$ws.Range("J2").Value = 1.073494826267255
$ws.Range("K2").Value = 1.074114844067272
$ws.Range("L2").Value = 1.083971358608884
$ws.Range("M2").Value = 1.089160079707898
$ws.Range("N2").Value = 1.028032011281404
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.069702091565755
$ws.Range("D3").Value = 1.07233618452538
$ws.Range("E3").Value = 1.08240465979342
$ws.Range("F3").Value = 1.087613104602405
$ws.Range("I3").Value = 1.058716993473628
$ws.Range("J3").Value = 1.074297180621242
$ws.Range("K3").Value = 1.07484970978311
$ws.Range("L3").Value = 1.084893568893259
$ws.Range("M3").Value = 1.090089471891792
$ws.Range("N3").Value = 1.028311167035495
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.070443957061496
$ws.Range("D4").Value = 1.072930529969989
$ws.Range("E4").Value = 1.083119682569805
$ws.Range("F4").Value = 1.088332183736879
$ws.Range("I4").Value = 1.058966535584599
$ws.Range("J4").Value = 1.074815785774452
$ws.Range("K4").Value = 1.075324553912134
$ws.Range("L4").Value = 1.085490049623539
$ws.Range("M4").Value = 1.090690632017859
$ws.Range("N4").Value = 1.02849135754835
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.070755842870683
$ws.Range("D5").Value = 1.073180385776736
$ws.Range("E5").Value = 1.083420361726843
$ws.Range("F5").Value = 1.088634575864414
$ws.Range("I5").Value = 1.059071132810084
$ws.Range("J5").Value = 1.075033670852813
$ws.Range("K5").Value = 1.075524019870868
$ws.Range("L5").Value = 1.085740750387119
$ws.Range("M5").Value = 1.090943307679589
$ws.Range("N5").Value = 1.028567003833789
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.070808210221698
$ws.Range("D6").Value = 1.073222337257383
$ws.Range("E6").Value = 1.083470851992206
$ws.Range("F6").Value = 1.088685354194145
$ws.Range("I6").Value = 1.059088676955128
$ws.Range("J6").Value = 1.075070246720483
$ws.Range("K6").Value = 1.075557501794802
$ws.Range("L6").Value = 1.085782840655843
$ws.Range("M6").Value = 1.090985729991751
$ws.Range("N6").Value = 1.028579698973237
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.070448124472006
$ws.Range("D7").Value = 1.07293386858296
$ws.Range("E7").Value = 1.083123699930536
$ws.Range("F7").Value = 1.088336223956288
$ws.Range("I7").Value = 1.058967934436787
$ws.Range("J7").Value = 1.074818697701547
$ws.Range("K7").Value = 1.075327219808579
$ws.Range("L7").Value = 1.085493399733684
$ws.Range("M7").Value = 1.090694008486556
$ws.Range("N7").Value = 1.02849236875347
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.068943069677355
$ws.Range("D8").Value = 1.071728056844762
$ws.Range("E8").Value = 1.081673369414488
$ws.Range("F8").Value = 1.086877690248778
$ws.Range("I8").Value = 1.058460616681324
$ws.Range("J8").Value = 1.073766104278329
$ws.Range("K8").Value = 1.074363332826638
$ws.Range("L8").Value = 1.0842830760731
$ws.Range("M8").Value = 1.089474217634984
$ws.Range("N8").Value = 1.028126444774473
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.066291081114362
$ws.Range("D9").Value = 1.069603030098203
$ws.Range("E9").Value = 1.079120347491209
$ws.Range("F9").Value = 1.084310461831131
$ws.Range("I9").Value = 1.057556717197595
$ws.Range("J9").Value = 1.071906904967158
$ws.Range("K9").Value = 1.072659750889959
$ws.Range("L9").Value = 1.082148392720213
$ws.Range("M9").Value = 1.087323100243855
$ws.Range("N9").Value = 1.027478255521939
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.064523105227756
$ws.Range("D10").Value = 1.068186179784113
$ws.Range("E10").Value = 1.077420085607603
$ws.Range("F10").Value = 1.082600898429942
$ws.Range("I10").Value = 1.056947420470979
$ws.Range("J10").Value = 1.070664458722647
$ws.Range("K10").Value = 1.071520589528765
$ws.Range("L10").Value = 1.080723942044013
$ws.Range("M10").Value = 1.085887862632484
$ws.Range("N10").Value = 1.027043853050268
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.063757542420853
$ws.Range("D11").Value = 1.067572625602943
$ws.Range("E11").Value = 1.076684262434334
$ws.Range("F11").Value = 1.081861088301438
$ws.Range("I11").Value = 1.056681996332068
$ws.Range("J11").Value = 1.070125752661171
$ws.Range("K11").Value = 1.071026499427752
$ws.Range("L11").Value = 1.08010681766745
$ws.Range("M11").Value = 1.085266107952095
$ws.Range("N11").Value = 1.026855211563699
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.063473174387989
$ws.Range("D12").Value = 1.067344716588287
$ws.Range("E12").Value = 1.076411004571024
$ws.Range("F12").Value = 1.081586355657146
$ws.Range("I12").Value = 1.056583166041367
$ws.Range("J12").Value = 1.069925544638085
$ws.Range("K12").Value = 1.070842847809659
$ws.Range("L12").Value = 1.079877540288561
$ws.Range("M12").Value = 1.085035116796233
$ws.Range("N12").Value = 1.026785060146459
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.063534172494228
$ws.Range("D13").Value = 1.067393604198684
$ws.Range("E13").Value = 1.076469616623126
$ws.Range("F13").Value = 1.081645283775716
$ws.Range("I13").Value = 1.056604376338289
$ws.Range("J13").Value = 1.069968494863552
$ws.Range("K13").Value = 1.070882247349583
$ws.Range("L13").Value = 1.079926723333578
$ws.Range("M13").Value = 1.085084667174497
$ws.Range("N13").Value = 1.026800111556209
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.063734036557891
$ws.Range("D14").Value = 1.067553786731338
$ws.Range("E14").Value = 1.076661673641727
$ws.Range("F14").Value = 1.081838377477238
$ws.Range("I14").Value = 1.056673831886406
$ws.Range("J14").Value = 1.070109205622744
$ws.Range("K14").Value = 1.071011321271492
$ws.Range("L14").Value = 1.080087866549676
$ws.Range("M14").Value = 1.085247015047835
$ws.Range("N14").Value = 1.026849414490174
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.063857178846271
$ws.Range("D15").Value = 1.067652479448696
$ws.Range("E15").Value = 1.076780014190045
$ws.Range("F15").Value = 1.081957357576548
$ws.Range("I15").Value = 1.05671659393547
$ws.Range("J15").Value = 1.070195887755099
$ws.Range("K15").Value = 1.071090831455013
$ws.Range("L15").Value = 1.080187145570659
$ws.Range("M15").Value = 1.085347037110525
$ws.Range("N15").Value = 1.026879780843569
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.064573912630347
$ws.Range("D16").Value = 1.068226898298391
$ws.Range("E16").Value = 1.077468928167563
$ws.Range("F16").Value = 1.082650006451351
$ws.Range("I16").Value = 1.056965002168034
$ws.Range("J16").Value = 1.070700195674332
$ws.Range("K16").Value = 1.071553363221892
$ws.Range("L16").Value = 1.080764891593454
$ws.Range("M16").Value = 1.085929120336677
$ws.Range("N16").Value = 1.027056361130255
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.065023495187561
$ws.Range("D17").Value = 1.068587202785985
$ws.Range("E17").Value = 1.077901172697917
$ws.Range("F17").Value = 1.083084604700854
$ws.Range("I17").Value = 1.057120394897615
$ws.Range("J17").Value = 1.071016341857559
$ws.Range("K17").Value = 1.071843275758537
$ws.Range("L17").Value = 1.081127208177176
$ws.Range("M17").Value = 1.086294168405971
$ws.Range("N17").Value = 1.027166980014369
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.065285727476178
$ws.Range("D18").Value = 1.068797357628687
$ws.Range("E18").Value = 1.078153332398886
$ws.Range("F18").Value = 1.08333814120695
$ws.Range("I18").Value = 1.057210879013371
$ws.Range("J18").Value = 1.071200675282047
$ws.Range("K18").Value = 1.072012297200355
$ws.Range("L18").Value = 1.081338509723226
$ws.Range("M18").Value = 1.086507067087469
$ws.Range("N18").Value = 1.027231449814958
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.065375141675618
$ws.Range("D19").Value = 1.068869014218291
$ws.Range("E19").Value = 1.078239318939063
$ws.Range("F19").Value = 1.083424597865798
$ws.Range("I19").Value = 1.057241705676871
$ws.Range("J19").Value = 1.071263516481926
$ws.Range("K19").Value = 1.072069915651727
$ws.Range("L19").Value = 1.081410552752068
$ws.Range("M19").Value = 1.08657965535962
$ws.Range("N19").Value = 1.0272534234787
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.064975259409611
$ws.Range("D20").Value = 1.068548546009795
$ws.Range("E20").Value = 1.077854792944188
$ws.Range("F20").Value = 1.083037972009507
$ws.Range("I20").Value = 1.057103738642763
$ws.Range("J20").Value = 1.0709824295361
$ws.Range("K20").Value = 1.071812179127663
$ws.Range("L20").Value = 1.081088338308475
$ws.Range("M20").Value = 1.086255005070418
$ws.Range("N20").Value = 1.027155117067594
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.063675181691082
$ws.Range("D21").Value = 1.067506617182669
$ws.Range("E21").Value = 1.076605115996192
$ws.Range("F21").Value = 1.081781514384262
$ws.Range("I21").Value = 1.05665338558532
$ws.Range("J21").Value = 1.070067772814466
$ws.Range("K21").Value = 1.070973315660926
$ws.Range("L21").Value = 1.080040415271882
$ws.Range("M21").Value = 1.085199208859152
$ws.Range("N21").Value = 1.026834898257397
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.062857747920232
$ws.Range("D22").Value = 1.066851470598916
$ws.Range("E22").Value = 1.075819739136427
$ws.Range("F22").Value = 1.080991909856445
$ws.Range("I22").Value = 1.056368842298756
$ws.Range("J22").Value = 1.069492063155038
$ws.Range("K22").Value = 1.070445168019527
$ws.Range("L22").Value = 1.079381255391337
$ws.Range("M22").Value = 1.084535134251594
$ws.Range("N22").Value = 1.026633092000631
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.063291087664353
$ws.Range("D23").Value = 1.067198780546087
$ws.Range("E23").Value = 1.07623604983116
$ws.Range("F23").Value = 1.08141045832815
$ws.Range("I23").Value = 1.056519815810904
$ws.Range("J23").Value = 1.069797317542544
$ws.Range("K23").Value = 1.070725217575894
$ws.Range("L23").Value = 1.0797307161352
$ws.Range("M23").Value = 1.084887197012356
$ws.Range("N23").Value = 1.026740118083318
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.064997055085974
$ws.Range("D24").Value = 1.068566013357224
$ws.Range("E24").Value = 1.077875749838415
$ws.Range("F24").Value = 1.083059043183249
$ws.Range("I24").Value = 1.057111265362543
$ws.Range("J24").Value = 1.070997753268299
$ws.Range("K24").Value = 1.071826230602468
$ws.Range("L24").Value = 1.081105902026354
$ws.Range("M24").Value = 1.086272701381273
$ws.Range("N24").Value = 1.027160477584214
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.066976676602372
$ws.Range("D25").Value = 1.070152428970521
$ws.Range("E25").Value = 1.07978005400969
$ws.Range("F25").Value = 1.084973811368787
$ws.Range("I25").Value = 1.057791576303986
$ws.Range("J25").Value = 1.072388075980418
$ws.Range("K25").Value = 1.07310077268085
$ws.Range("L25").Value = 1.082700491840345
$ws.Range("M25").Value = 1.087879418307189
$ws.Range("N25").Value = 1.027646229137347
